# Edit script for wymagania-projekt-hih.xlsx
# Implements:
#  - A2 rich text: last run "OK?" -> "OK" (keeping all run-level formatting)
#  - C2 (Czy spelnione for "strona responsywna..."): 0 -> 1
#  - C11 (Czy spelnione for "animacja wykorzystujaca..."): 0 -> 1
#  - Active cell selection: C5 -> C3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the rich text in A2: change trailing "OK?" to "OK" ---
$cell = $ws.Range("A2")

$newFullText = "strona responsywna (minimum to obsługa dwóch różnych wielkości ekranu) (1,6pkt). Minimalne wymagania na responsywność:`no viewport, - OK`no media queries (@media) - OK`no zmiana położenia menu - PRZEMYŚLEĆ`no wykorzystanie viewport-width property - OK"
$cell.Characters().Text = $newFullText

$colorGreen  = 5287936   # RGB(0,176,80)   -> FF00B050
$colorOrange = 49407     # RGB(255,192,0)  -> FFFFC000
$colorBlack  = 0         # RGB(0,0,0)      -> theme text 1 / black

# NOTE: Characters(1, 133) (the leading run) keeps the default cell font -
# do not touch it so it stays without an explicit <rPr>, matching the
# original formatting of that run.
$cell.Characters(134, 2).Font.Color   = $colorGreen
$cell.Characters(136, 28).Font.Color  = $colorBlack
$cell.Characters(164, 2).Font.Color   = $colorGreen
$cell.Characters(166, 27).Font.Color  = $colorBlack
$cell.Characters(193, 10).Font.Color  = $colorOrange
$cell.Characters(203, 43).Font.Color  = $colorBlack
$cell.Characters(246, 2).Font.Color   = $colorOrange

# --- 2. Mark requirement in row 2 as fulfilled ---
$ws.Range("C2").Value = 1

# --- 3. Mark requirement in row 11 as fulfilled ---
$ws.Range("C11").Value = 1

# Recalculate dependent formulas (Uzyskane punkty column + Suma row)
$wb.Application.Calculate()

# --- 4. Update the active selection shown when the sheet was last saved ---
$ws.Range("C3").Select()
